# This script updates the "Organization" column (A) for the ranked list of
# funders/sponsors in the active worksheet, reflecting a refreshed
# similarity/embedding-based ranking (rank_group values in column B are
# unchanged; only the organization names assigned to each rank position
# change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Owl Ventures'
$ws.Range("A3").Value = 'Bill & Melinda Gates Foundation'
$ws.Range("A4").Value = 'Y Combinator'
$ws.Range("A5").Value = 'Reach Capital'
$ws.Range("A6").Value = 'Google for Education'
$ws.Range("A7").Value = 'U.S. Dept. of Education (EIR Program)'
$ws.Range("A8").Value = 'National Science Foundation (NSF)'
$ws.Range("A9").Value = '500 Global Flagship VC (non-accelerator checks)'
$ws.Range("A10").Value = 'IES SBIR (ED/IES)'
$ws.Range("A11").Value = '500 Global (seed/accelerator)'
$ws.Range("A12").Value = 'Buffalo Sabres Foundation'
$ws.Range("A13").Value = 'TGR Foundation (Tiger Woods)'
$ws.Range("A14").Value = 'Berkeley SkyDeck Fund (UC Berkeley)'
$ws.Range("A15").Value = 'Penn State University - Outreach & Engagement'
$ws.Range("A16").Value = 'Chan Zuckerberg Initiative (CZI)'
$ws.Range("A17").Value = 'Nashville Predators Foundation'
$ws.Range("A18").Value = 'Portland Trail Blazers Foundation'
$ws.Range("A19").Value = 'NFL Foundation'
$ws.Range("A20").Value = 'San Jose Sharks Foundation'
$ws.Range("A21").Value = 'San Antonio Spurs - Spurs Give'
$ws.Range("A22").Value = 'Indiana Pacers Foundation'
$ws.Range("A23").Value = 'Cleveland Cavaliers Community Foundation'
$ws.Range("A24").Value = 'Austin FC - 4ATX Foundation'
$ws.Range("A25").Value = 'Baltimore Ravens Foundation'
$ws.Range("A26").Value = 'Toronto FC - MLSE Foundation'
$ws.Range("A27").Value = 'Houston Texans Foundation'
$ws.Range("A28").Value = 'Philadelphia Eagles Foundation'
$ws.Range("A29").Value = 'Oakland Roots SC'
$ws.Range("A30").Value = 'Jacksonville Jaguars Foundation'
$ws.Range("A31").Value = 'Florida State University Research Foundation'
